$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-13 and add new rows 14-17
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adam15"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 37.63904266666666
$ws.Range("H2").Value = 112.917128
$ws.Range("I2").Value = 0.4850220755088102
$ws.Range("J2").Value = 0.4850220755088102
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.481489333333333
$ws.Range("N2").Value = 7.444467999999999
$ws.Range("O2").Value = 0.2345069082418988
$ws.Range("P2").Value = 0.2345069082418987
$ws.Range("Q2").Value = 93.40088289421152
$ws.Range("R2").Value = 840.6079460479037
$ws.Range("S2").Value = 0.1137410273566398
$ws.Range("T2").Value = 0.1137410273566398

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adam15"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 37.63904266666666
$ws.Range("H3").Value = 112.917128
$ws.Range("I3").Value = 0.4850220755088102
$ws.Range("J3").Value = 0.4850220755088102
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.245227
$ws.Range("N3").Value = 21.735681
$ws.Range("O3").Value = 0.6846919551326144
$ws.Range("P3").Value = 0.6846919551326142
$ws.Range("Q3").Value = 272.7034081826853
$ws.Range("R3").Value = 2454.330673644168
$ws.Range("S3").Value = 0.3320907131626057
$ws.Range("T3").Value = 0.3320907131626057

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adam15"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 37.63904266666666
$ws.Range("H4").Value = 112.917128
$ws.Range("I4").Value = 0.4850220755088102
$ws.Range("J4").Value = 0.4850220755088102
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2001876666666667
$ws.Range("N4").Value = 0.600563
$ws.Range("O4").Value = 0.01891823194544989
$ws.Range("P4").Value = 0.01891823194544989
$ws.Range("Q4").Value = 7.53487212700711
$ws.Range("R4").Value = 67.81384914306399
$ws.Range("S4").Value = 0.009175760123139183
$ws.Range("T4").Value = 0.009175760123139181

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Adam15"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 37.63904266666666
$ws.Range("H5").Value = 112.917128
$ws.Range("I5").Value = 0.4850220755088102
$ws.Range("J5").Value = 0.4850220755088102
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6548283333333333
$ws.Range("N5").Value = 1.964485
$ws.Range("O5").Value = 0.06188290468003712
$ws.Range("P5").Value = 0.06188290468003711
$ws.Range("Q5").Value = 24.64711157767555
$ws.Range("R5").Value = 221.82400419908
$ws.Range("S5").Value = 0.03001457486642546
$ws.Range("T5").Value = 0.03001457486642546

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Adam15"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.57434666666667
$ws.Range("H6").Value = 52.72304
$ws.Range("I6").Value = 0.2264655392929762
$ws.Range("J6").Value = 0.2264655392929762
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.481489333333333
$ws.Range("N6").Value = 7.444467999999999
$ws.Range("O6").Value = 0.2345069082418988
$ws.Range("P6").Value = 0.2345069082418987
$ws.Range("Q6").Value = 43.61055379363555
$ws.Range("R6").Value = 392.4949841427199
$ws.Range("S6").Value = 0.05310773344293009
$ws.Range("T6").Value = 0.05310773344293007

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Adam15"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.57434666666667
$ws.Range("H7").Value = 52.72304
$ws.Range("I7").Value = 0.2264655392929762
$ws.Range("J7").Value = 0.2264655392929762
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.245227
$ws.Range("N7").Value = 21.735681
$ws.Range("O7").Value = 0.6846919551326144
$ws.Range("P7").Value = 0.6846919551326142
$ws.Range("Q7").Value = 127.3301309766933
$ws.Range("R7").Value = 1145.97117879024
$ws.Range("S7").Value = 0.1550591328686698
$ws.Range("T7").Value = 0.1550591328686697

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Adam15"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 17.57434666666667
$ws.Range("H8").Value = 52.72304
$ws.Range("I8").Value = 0.2264655392929762
$ws.Range("J8").Value = 0.2264655392929762
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2001876666666667
$ws.Range("N8").Value = 0.600563
$ws.Range("O8").Value = 0.01891823194544989
$ws.Range("P8").Value = 0.01891823194544989
$ws.Range("Q8").Value = 3.518167452391111
$ws.Range("R8").Value = 31.66350707152
$ws.Range("S8").Value = 0.00428432759999592
$ws.Range("T8").Value = 0.004284327599995919

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Adam15"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 17.57434666666667
$ws.Range("H9").Value = 52.72304
$ws.Range("I9").Value = 0.2264655392929762
$ws.Range("J9").Value = 0.2264655392929762
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6548283333333333
$ws.Range("N9").Value = 1.964485
$ws.Range("O9").Value = 0.06188290468003712
$ws.Range("P9").Value = 0.06188290468003711
$ws.Range("Q9").Value = 11.50818013715556
$ws.Range("R9").Value = 103.5736212344
$ws.Range("S9").Value = 0.01401434538138045
$ws.Range("T9").Value = 0.01401434538138044

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Adam15"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 19.168158
$ws.Range("H10").Value = 57.504474
$ws.Range("I10").Value = 0.2470036195972184
$ws.Range("J10").Value = 0.2470036195972184
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.481489333333333
$ws.Range("N10").Value = 7.444467999999999
$ws.Range("O10").Value = 0.2345069082418988
$ws.Range("P10").Value = 0.2345069082418987
$ws.Range("Q10").Value = 47.56557961664799
$ws.Range("R10").Value = 428.090216549832
$ws.Range("S10").Value = 0.05792405515630176
$ws.Range("T10").Value = 0.05792405515630174

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Adam15"
$ws.Range("C11").Value = "Itgb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 19.168158
$ws.Range("H11").Value = 57.504474
$ws.Range("I11").Value = 0.2470036195972184
$ws.Range("J11").Value = 0.2470036195972184
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.245227
$ws.Range("N11").Value = 21.735681
$ws.Range("O11").Value = 0.6846919551326144
$ws.Range("P11").Value = 0.6846919551326142
$ws.Range("Q11").Value = 138.877655881866
$ws.Range("R11").Value = 1249.898902936794
$ws.Range("S11").Value = 0.169121391226852
$ws.Range("T11").Value = 0.169121391226852

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Adam15"
$ws.Range("C12").Value = "Itgb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 19.168158
$ws.Range("H12").Value = 57.504474
$ws.Range("I12").Value = 0.2470036195972184
$ws.Range("J12").Value = 0.2470036195972184
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2001876666666667
$ws.Range("N12").Value = 0.600563
$ws.Range("O12").Value = 0.01891823194544989
$ws.Range("P12").Value = 0.01891823194544989
$ws.Range("Q12").Value = 3.837228824318
$ws.Range("R12").Value = 34.535059418862
$ws.Range("S12").Value = 0.00467287176690585
$ws.Range("T12").Value = 0.004672871766905849

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Adam15"
$ws.Range("C13").Value = "Itgb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 19.168158
$ws.Range("H13").Value = 57.504474
$ws.Range("I13").Value = 0.2470036195972184
$ws.Range("J13").Value = 0.2470036195972184
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6548283333333333
$ws.Range("N13").Value = 1.964485
$ws.Range("O13").Value = 0.06188290468003712
$ws.Range("P13").Value = 0.06188290468003711
$ws.Range("Q13").Value = 12.55185295621
$ws.Range("R13").Value = 112.96667660589
$ws.Range("S13").Value = 0.01528530144715881
$ws.Range("T13").Value = 0.01528530144715881

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Adam15"
$ws.Range("C14").Value = "Itgb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3.221194
$ws.Range("H14").Value = 9.663582
$ws.Range("I14").Value = 0.04150876560099527
$ws.Range("J14").Value = 0.04150876560099527
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.481489333333333
$ws.Range("N14").Value = 7.444467999999999
$ws.Range("O14").Value = 0.2345069082418988
$ws.Range("P14").Value = 0.2345069082418987
$ws.Range("Q14").Value = 7.993358551597332
$ws.Range("R14").Value = 71.94022696437598
$ws.Range("S14").Value = 0.009734092286027082
$ws.Range("T14").Value = 0.00973409228602708

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Adam15"
$ws.Range("C15").Value = "Itgb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.221194
$ws.Range("H15").Value = 9.663582
$ws.Range("I15").Value = 0.04150876560099527
$ws.Range("J15").Value = 0.04150876560099527
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 7.245227
$ws.Range("N15").Value = 21.735681
$ws.Range("O15").Value = 0.6846919551326144
$ws.Range("P15").Value = 0.6846919551326142
$ws.Range("Q15").Value = 23.338281741038
$ws.Range("R15").Value = 210.044535669342
$ws.Range("S15").Value = 0.02842071787448686
$ws.Range("T15").Value = 0.02842071787448685

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Adam15"
$ws.Range("C16").Value = "Itgb3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 3.221194
$ws.Range("H16").Value = 9.663582
$ws.Range("I16").Value = 0.04150876560099527
$ws.Range("J16").Value = 0.04150876560099527
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.2001876666666667
$ws.Range("N16").Value = 0.600563
$ws.Range("O16").Value = 0.01891823194544989
$ws.Range("P16").Value = 0.01891823194544989
$ws.Range("Q16").Value = 0.6448433107406667
$ws.Range("R16").Value = 5.803589796666
$ws.Range("S16").Value = 0.0007852724554089403
$ws.Range("T16").Value = 0.0007852724554089402

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Adam15"
$ws.Range("C17").Value = "Itgb3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 3.221194
$ws.Range("H17").Value = 9.663582
$ws.Range("I17").Value = 0.04150876560099527
$ws.Range("J17").Value = 0.04150876560099527
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6548283333333333
$ws.Range("N17").Value = 1.964485
$ws.Range("O17").Value = 0.06188290468003712
$ws.Range("P17").Value = 0.06188290468003711
$ws.Range("Q17").Value = 2.109329098363333
$ws.Range("R17").Value = 18.98396188527
$ws.Range("S17").Value = 0.002568682985072394
$ws.Range("T17").Value = 0.002568682985072394
